# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.Formula = "'55.910.17"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +8.49%  "

# Row 3
$c = $ws.Range("D3")
$c.Formula = "'3.218.93"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.69%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$c = $ws.Range("D5")
$c.Formula = "'398.25"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.30%  "

# Row 6
$c = $ws.Range("D6")
$c.Formula = "'109.73"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.20%  "

# Row 7
$ws.Range("E7").Value = "  +2.91%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$c = $ws.Range("D9")
$c.Formula = "'0.622"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.42%  "

# Row 10
$c = $ws.Range("D10")
$c.Formula = "'39.27"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.89%  "

# Row 11
$c = $ws.Range("D11")
$c.Formula = "'0.0900"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.36%  "

# Row 13
$c = $ws.Range("D13")
$c.Formula = "'3.724.87"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.66%  "

# Row 14
$c = $ws.Range("D14")
$c.Formula = "'19.07"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.52%  "

# Row 15
$ws.Range("E15").Value = "  +2.96%  "

# Row 16
$c = $ws.Range("D16")
$c.Formula = "'3.218.90"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.82%  "

# Row 17
$ws.Range("E17").Value = "  +6.05%  "

# Row 18
$c = $ws.Range("D18")
$c.Formula = "'10.57"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.71%  "

# Row 19
$c = $ws.Range("D19")
$c.Formula = "'55.787.92"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +8.21%  "

# Row 20
$c = $ws.Range("D20")
$c.Formula = "'3.36"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.50%  "

# Row 21
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D21")
$c.Formula = "'13.09"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.84%  "

# Row 22
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D22")
$c.Formula = "'0.0000102"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +5.94%  "

# Row 23
$c = $ws.Range("D23")
$c.Formula = "'303.47"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +13.96%  "

# Row 24
$c = $ws.Range("D24")
$c.Formula = "'75.05"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +7.38%  "

# Row 25
$ws.Range("E25").Value = "  +2.18%  "

# Row 26
$c = $ws.Range("D26")
$c.Formula = "'8.25"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.76%  "

# Row 27
$c = $ws.Range("D27")
$c.Formula = "'28.26"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.49%  "

# Row 28
$c = $ws.Range("D28")
$c.Formula = "'7.52"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.55%  "

# Row 29
$c = $ws.Range("D29")
$c.Formula = "'0.173"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.86%  "

# Row 30
$ws.Range("E30").Value = "  -0.20%  "

# Row 31
$ws.Range("E31").Value = "  +9.95%  "

# Row 32
$ws.Range("E32").Value = "  +3.11%  "

# Row 33
$c = $ws.Range("D33")
$c.Formula = "'0.0494"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.87%  "

# Row 34
$c = $ws.Range("D34")
$c.Formula = "'36.20"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.86%  "

# Row 35
$ws.Range("E35").Value = "  +2.28%  "

# Row 36
$c = $ws.Range("D36")
$c.Formula = "'51.40"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.32%  "

# Row 37
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D37")
$c.Formula = "'3.10"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +23.34%  "

# Row 38
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D38")
$c.Formula = "'1.00"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D39")
$c.Formula = "'3.51"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.20%  "

# Row 40
$c = $ws.Range("D40")
$c.Formula = "'134.99"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.62%  "

# Row 41
$c = $ws.Range("D41")
$c.Formula = "'4.04"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +10.07%  "

# Row 42
$ws.Range("E42").Value = "  +1.79%  "

# Row 43
$ws.Range("E43").Value = "  +2.92%  "

# Row 44
$c = $ws.Range("D44")
$c.Formula = "'17.04"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.01%  "

# Row 45
$c = $ws.Range("D45")
$c.Formula = "'0.284"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.94%  "

# Row 46
$c = $ws.Range("D46")
$c.Formula = "'22.25"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.98%  "

# Row 47
$c = $ws.Range("D47")
$c.Formula = "'2.11"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.16%  "

# Row 48
$c = $ws.Range("D48")
$c.Formula = "'2.151.21"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.02%  "

# Row 49
$c = $ws.Range("D49")
$c.Formula = "'2.47"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "

# Row 50
$c = $ws.Range("D50")
$c.Formula = "'2.09"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +43.39%  "

# Row 51
$c = $ws.Range("D51")
$c.Formula = "'0.0361"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +8.89%  "
